$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows down (bottom-up, so we don't overwrite data before it's moved):
# old A3 (Ranar the Ever-Watchful) -> A4
$ws.Range("A4").Value = $ws.Range("A3").Text
# old A2 (Lathril, Blade of the Elves) -> A3
$ws.Range("A3").Value = $ws.Range("A2").Text

# New row 2: Inspired Sphinx
$ws.Range("A2").Value = "('Inspired Sphinx', ['{5}{U}{U}', 'Creature — Sphinx', 'Flying', 'When Inspired Sphinx enters the battlefield, draw cards equal to the number of opponents you have.', '{3}{U}: Create a 1/1 colorless Thopter artifact creature token with flying.', '5/5'])"

# New row 5 (appended at the end): Wolverine Riders
$ws.Range("A5").Value = "('Wolverine Riders', ['{4}{G}{G}', 'Creature — Elf Warrior', 'At the beginning of each upkeep, create a 1/1 green Elf Warrior creature token.', 'Whenever another Elf enters the battlefield under your control, you gain life equal to its toughness.', '4/4'])"
